# "Done VNP into addPayment" - replace renter/vehicle sample data
# in the contract body. Each edit targets the exact character Range of
# the original run(s) (offsets computed from the pristine document),
# applied highest-offset-first so earlier offsets stay valid as the
# document length changes.

$d = $word.ActiveDocument

# Rental price part C (bold currency suffix)
$r = $d.Range(1253, 1262)
$r.Text = "NĐ/ngày."

# Rental price part B
$r = $d.Range(1247, 1253)
$r.Text = "………….V"

# Rental price part A: 1187950.00 -> 944150.00
$r = $d.Range(1237, 1247)
$r.Text = "944150.00…"

# Plate part B: 51.. -> 32..
$r = $d.Range(1143, 1147)
$r.Text = "32……"

# Plate part A: 1K-392 -> 1L-288
$r = $d.Range(1136, 1142)
$r.Text = "1L-288"

# Vehicle year last digit: 2 -> 0 (bold)
$r = $d.Range(1113, 1114)
$r.Text = "0"

# Vehicle brand part B: reta 202 -> nova 202
$r = $d.Range(1105, 1113)
$r.Text = "nova 202"

# Vehicle brand part A: Hyundai C -> Toyota In
$r = $d.Range(1096, 1105)
$r.Text = "Toyota In"

# Renter 2 phone: 0777444499 -> 0999998888
$r = $d.Range(736, 757)
$r.Text = "0999998888……………………………"

# Renter 2 name part C
$r = $d.Range(692, 694)
$r.Text = "…"

# Renter 2 name part B
$r = $d.Range(674, 692)
$r.Text = "…………Giới tính:…………"

# Renter 2 name part A: Phan Van Khai -> Van Huu Toan
$r = $d.Range(660, 674)
$r.Text = "Van Huu Toan……"

# Renter 1 phone: 0969998878 -> 0999777712
$r = $d.Range(595, 613)
$r.Text = "0999777712……………………"

# Renter 1 name part C
$r = $d.Range(466, 476)
$r.Text = "…………."

# Renter 1 name part B
$r = $d.Range(447, 466)
$r.Text = "i…………………Giới tính:…"

# Renter 1 name part A: Nguyen Pham -> Quang La Tu
$r = $d.Range(436, 447)
$r.Text = "Quang La Tu"
